$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "318.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.59%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.210"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.76%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08237"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.89%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.153"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.21%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.063"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.44%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9271"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.83%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1021"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.14%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1888"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.85%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09158"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.05%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03619"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.27%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09915"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.37%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001437"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.71%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005682"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.01%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.466"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.07%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.135"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.67%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.799"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17.00%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.61%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1300"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.42%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.062"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.88%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2189"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.61%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04596"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.34%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.91%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.73%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004501"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.29%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02008"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.65%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04960"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.76%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007811"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.13%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.01%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007586"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.93%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002097"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01192"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.42%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.09%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.04"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-14.33%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.02%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.01%"
